$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row (row 11) appended to the results table
$ws.Range("A11").Value = 46059
$ws.Range("A11").NumberFormat = "m/d/yy"

$ws.Range("B11").Value = "Volta a la Comunitat Valenciana"
$ws.Range("C11").Value = "Stage 3"
$ws.Range("D11").Value = "Andrew August"
$ws.Range("E11").Value = "Ådne Holter"
$ws.Range("F11").Value = "Florian Vermeersch"
$ws.Range("G11").Value = "Jonathan Vervenne"
$ws.Range("H11").Value = "Biniam Girmay"
$ws.Range("I11").Value = "Ben Turner"
$ws.Range("J11").Value = "Magnus Cort"
$ws.Range("K11").Value = "Mathias Vacek"
$ws.Range("L11").Value = "Mirco Maestri"
$ws.Range("M11").Value = "Aleksandr Vlasov"

# Update the active selection to match the saved workbook state
$ws.Range("D24").Select()
